$d = $word.ActiveDocument

# 1) Insert new paragraph right after the paragraph ending with
#    "I BILAGA 1 finns artfakta om fridlysta arter." containing the
#    "Vi förväntar oss..." text.
$newText = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*I BILAGA 1 finns artfakta om fridlysta arter.*") {
        $insertionRange = $p.Range
        $insertionRange.Collapse(0)
        $insertionRange.InsertParagraphAfter()
        $newP = $p.Next()
        $newP.Range.Text = $newText
        break
    }
}

# 2) Remove the old "Vi förväntar oss..." paragraph (near end of doc,
#    right before the page-break paragraph), plus the two empty
#    paragraphs preceding it. Distinguish it from the freshly inserted
#    copy by checking that the paragraph right before it is empty.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Vi förväntar oss att ni återkommer*") {
        $prev1 = $p.Previous()
        if ($prev1.Range.Text.Trim() -eq "") {
            $prev2 = $prev1.Previous()
            $startRange = $prev2.Range
            $endRange = $p.Range
            $delRange = $d.Range($startRange.Start, $endRange.End)
            $delRange.Delete()
            break
        }
    }
}

# 3) Update the date in the header from 2023-11-13 to 2023-11-14.
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "2023-11-14", 2) | Out-Null
        }
    }
}
